$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 260140.53
$ws.Range("I132").Value = 3850.6758
$ws.Range("K132").Value = 11552.0274
$ws.Range("M132").Value = -9022.027399999999

$ws.Range("H137").Value = 48850.047
$ws.Range("I137").Value = 92085.27
$ws.Range("J137").Value = 1291.3
$ws.Range("K137").Value = 276255.81
$ws.Range("L137").Value = 3873.9
$ws.Range("M137").Value = -273705.81
$ws.Range("N137").Value = -8973.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = 0

$ws.Range("H5").Value = 7166
$ws.Range("I5").Value = 9999
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 9999
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = -9887
$ws.Range("N5").Value = -1724

$ws.Range("H9").Value = 12839.4
$ws.Range("J9").Value = 12839.4
$ws.Range("L9").Value = 12839.4
$ws.Range("N9").Value = -13179.4

$ws.Range("H20").Value = 12839.4
$ws.Range("J20").Value = 12839.4
$ws.Range("L20").Value = 12839.4
$ws.Range("N20").Value = -13379.4

$ws.Range("H23").Value = 80006.5
$ws.Range("I23").Value = 70006
$ws.Range("J23").Value = 90007
$ws.Range("K23").Value = 70006
$ws.Range("L23").Value = 90007
$ws.Range("M23").Value = -69747
$ws.Range("N23").Value = -90525

$ws.Range("H32").Value = 7404.2856
$ws.Range("I32").Value = 4158.031
$ws.Range("J32").Value = 49605.6
$ws.Range("K32").Value = 4158.031
$ws.Range("L32").Value = 49605.6
$ws.Range("M32").Value = -3871.031
$ws.Range("N32").Value = -50179.6

$ws.Range("H63").Value = 2110.9167
$ws.Range("I63").Value = 2110.9167
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2110.9167
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = ""
$ws.Range("N63").Value = -1424.9167

$ws.Range("H66").Value = 2110.9167
$ws.Range("I66").Value = 2110.9167
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10554.5835
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -7122.583500000001
$ws.Range("N66").Value = ""

$ws.Range("H74").Value = 61729.234
$ws.Range("I74").Value = 69615.8
$ws.Range("J74").Value = 2580
$ws.Range("K74").Value = 69615.8
$ws.Range("L74").Value = 2580
$ws.Range("M74").Value = -68741.8
$ws.Range("N74").Value = -4328

$ws.Range("H77").Value = 61729.234
$ws.Range("I77").Value = 69615.8
$ws.Range("J77").Value = 2580
$ws.Range("K77").Value = 348079
$ws.Range("L77").Value = 12900
$ws.Range("M77").Value = -343711
$ws.Range("N77").Value = -21636

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 7166
$ws.Range("I4").Value = 9999
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 9999
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -9884
$ws.Range("N4").Value = -1730

$ws.Range("H15").Value = 18995
$ws.Range("J15").Value = 18995
$ws.Range("L15").Value = 18995
$ws.Range("N15").Value = -19449

$ws.Range("H19").Value = 34000
$ws.Range("J19").Value = 34000
$ws.Range("L19").Value = 34000
$ws.Range("N19").Value = -34346

$ws.Range("H35").Value = 18516
$ws.Range("J35").Value = 18516
$ws.Range("L35").Value = 18516
$ws.Range("N35").Value = -19136

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 252122.75
$ws.Range("I6").Value = 500245.5
$ws.Range("J6").Value = 4000
$ws.Range("K6").Value = 500245.5
$ws.Range("L6").Value = 4000
$ws.Range("M6").Value = -500132.5
$ws.Range("N6").Value = -4226

$ws.Range("H7").Value = 571.4737
$ws.Range("I7").Value = 753.4286
$ws.Range("J7").Value = 62
$ws.Range("K7").Value = 753.4286
$ws.Range("L7").Value = 62
$ws.Range("M7").Value = -640.4286
$ws.Range("N7").Value = -288

$ws.Range("H31").Value = 8460.611999999999
$ws.Range("I31").Value = 9300.200000000001
$ws.Range("K31").Value = 9300.200000000001
$ws.Range("M31").Value = -9005.200000000001

$ws.Range("H34").Value = 8460.611999999999
$ws.Range("I34").Value = 9300.200000000001
$ws.Range("K34").Value = 9300.200000000001
$ws.Range("M34").Value = -9098.200000000001

$ws.Range("H62").Value = 1950
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 1950
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = ""
$ws.Range("M62").Value = 1950
$ws.Range("N62").Value = -3198

$ws.Range("H65").Value = 1950
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 1950
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = ""
$ws.Range("M65").Value = 9750
$ws.Range("N65").Value = -15990

$ws.Range("H107").Value = 588.7222
$ws.Range("I107").Value = 509.9
$ws.Range("J107").Value = 687.25
$ws.Range("K107").Value = 509.9
$ws.Range("L107").Value = 687.25
$ws.Range("M107").Value = 1410.1
$ws.Range("N107").Value = -4527.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12520
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 50600
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 50600
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -52596

$ws.Range("H83").Value = 12520
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 50600
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 253000
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -262984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 113700.445
$ws.Range("I40").Value = 2846
$ws.Range("K40").Value = 2846
$ws.Range("M40").Value = -2710

$ws.Range("H132").Value = 270212.94
$ws.Range("I132").Value = 79783.84
$ws.Range("J132").Value = 530800.1
$ws.Range("K132").Value = 239351.52
$ws.Range("L132").Value = 1592400.3
$ws.Range("M132").Value = -236821.52
$ws.Range("N132").Value = -1597460.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 181
$ws.Range("I107").Value = 173.33333
$ws.Range("J107").Value = 198.25
$ws.Range("K107").Value = 519.99999
$ws.Range("L107").Value = 594.75
$ws.Range("M107").Value = 1400.00001
$ws.Range("N107").Value = -4434.75
